$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rev B")

# Row 10 (J1 - Power Input Connector): no longer need to call out JST-RCY-Female;
# source is now Amazon.
$ws.Range("B10").ClearContents()
$ws.Range("F10").Value = "(Amazon)"

# Row 11 (J2 - Power Output Connector): the connector itself is supplied together
# with J1's cable, and the board-level cost tracking for it is dropped; the link
# now points at RapidLED, the actual supplier of the output connector.
$ws.Range("B11").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("F11").Value = "(RapidLED)"

# Row 12 (J3 - Dimming Pot Connector, onboard)
$ws.Range("F12").Value = "(Amazon)"

# Row 17 (J4 - Dimming Pot Connector, offboard)
$ws.Range("F17").Value = "(Amazon)"

# Row 19 (Board) - note the PCB fab source
$ws.Range("F19").Value = "(PCBway)"

$ws.Range("H21").Select()
